# Update edited session - rename the sheet and drop the oldest scan-log row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Pediatrics" -> "Session"
$ws.Name = "Session"

# Row 2 (Student ID 201404 / 14:29:44) is removed; the two remaining log
# rows (201447/14:30:05 and 201452/14:30:06) shift up into rows 2 and 3,
# and the used range shrinks from A1:F4 to A1:F3 accordingly.
$ws.Rows.Item(2).Delete()
